$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the station id for the Bowman, ND row (row 26): KBPP -> KBWW
$ws.Range("B26").Value = "KBWW"

# Update the Bowman Arpt, ND location (latitude/longitude)
$ws.Range("E26").Value = 46.1655
$ws.Range("F26").Value = -103.3

# The station-name cell (D26) picks up a new (near-duplicate) font/style,
# matching the theme minor font with explicit black color and vertical-center alignment.
$d26 = $ws.Range("D26")
$d26.Font.ThemeFont = 1
$d26.Font.Color = 0

# Move the active selection to E24, matching the saved selection state.
$ws.Range("E24").Select()
